$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 249, pushing the existing rows 249-251 down to 252-254
$ws.Rows("249:251").Insert()

# Fill in the new row 249 (Extra) with updated values for the new date/prices
$ws.Range("A249").Value = 5
$ws.Range("B249").Value = "Macroferia Regional de Talca"
$ws.Range("C249").Value = "Maule"
$ws.Range("D249").Value = "1/13/2023"
$ws.Range("E249").Value = 7
$ws.Range("F249").Value = 100112028
$ws.Range("G249").Value = "Sandia"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Extra"
$ws.Range("J249").Value = 1500
$ws.Range("K249").Value = 3000
$ws.Range("L249").Value = 3000
$ws.Range("M249").Value = 3000
$ws.Range("N249").Value = "`$/unidad"
$ws.Range("O249").Value = "Región del Maule"
$ws.Range("P249").Value = 3000
$ws.Range("Q249").Value = 1
$ws.Range("R249").Value = "Hortaliza"

# Fill in the new row 250 (Primera) with updated values for the new date/prices
$ws.Range("A250").Value = 5
$ws.Range("B250").Value = "Macroferia Regional de Talca"
$ws.Range("C250").Value = "Maule"
$ws.Range("D250").Value = "1/13/2023"
$ws.Range("E250").Value = 7
$ws.Range("F250").Value = 100112028
$ws.Range("G250").Value = "Sandia"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 4000
$ws.Range("K250").Value = 2300
$ws.Range("L250").Value = 2300
$ws.Range("M250").Value = 2300
$ws.Range("N250").Value = "`$/unidad"
$ws.Range("O250").Value = "Región del Maule"
$ws.Range("P250").Value = 2300
$ws.Range("Q250").Value = 1
$ws.Range("R250").Value = "Hortaliza"

# Fill in the new row 251 (Segunda) with updated values for the new date/prices
$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = "1/13/2023"
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112028
$ws.Range("G251").Value = "Sandia"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Segunda"
$ws.Range("J251").Value = 4000
$ws.Range("K251").Value = 800
$ws.Range("L251").Value = 800
$ws.Range("M251").Value = 800
$ws.Range("N251").Value = "`$/unidad"
$ws.Range("O251").Value = "Región del Maule"
$ws.Range("P251").Value = 800
$ws.Range("Q251").Value = 1
$ws.Range("R251").Value = "Hortaliza"
